$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "ChangePassword" worksheet after the last existing sheet
# ---------------------------------------------------------------------------
$sheetCount  = $wb.Worksheets.Count
$lastSheet   = $wb.Worksheets.Item($sheetCount)
$wsChange    = $wb.Worksheets.Add($null, $lastSheet)
$wsChange.Name = "ChangePassword"

# Header row (write in this order so the shared-string table is built
# New Password, Old Password, Confirm Password -> matches target layout)
$wsChange.Range("B1").Value = "New Password"
$wsChange.Range("A1").Value = "Old Password"
$wsChange.Range("C1").Value = "Confirm Password"

# Data row
$wsChange.Range("A2").Value = "test123"
$wsChange.Range("B2").Value = "test1234"
$wsChange.Range("C2").Value = "test1234"

# ---------------------------------------------------------------------------
# 2. Update the Login sheet's password value, re-using the "test123" string
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("B2").Value = "test123"

# ---------------------------------------------------------------------------
# 3. Apply the existing "Hyperlink" look (style index used by Login!A2/B2)
#    to the new password row, then drop the actual hyperlinks from Login
# ---------------------------------------------------------------------------
$null = $wsLogin.Range("A2").Copy()
$null = $wsChange.Range("A2:C2").PasteSpecial(-4122)

$null = $wsLogin.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 4. Cosmetic sheet formatting for ChangePassword (columns / print setup)
# ---------------------------------------------------------------------------
$wsChange.Columns.Item(1).ColumnWidth = 12.25
$wsChange.Columns.Item(2).ColumnWidth = 13.25
$wsChange.Columns.Item(3).ColumnWidth = 16.42

$wsChange.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Selections / active tab: Login -> B2, ChangePassword -> D2 (active)
# ---------------------------------------------------------------------------
$null = $wsLogin.Range("B2").Select()
$null = $wsChange.Range("D2").Select()
